$wb = $excel.ActiveWorkbook
$about = $wb.Worksheets.Item("About")
$bpp   = $wb.Worksheets.Item("BPP")

# Insert 5 new rows before the old "Notes:" row (row 8), shifting everything
# below it down by 5 rows (old row 8 -> 13, old row 9 -> 14, old row 11 -> 16).
# Excel automatically repoints the cross-sheet formulas that reference About!$A$11
# (on both the BPP and SYBPP sheets) to the shifted row.
$about.Rows("8:12").Insert()

# New second source citation block (mirrors the existing rows 3-6 pattern).
$about.Range("B8").Value = "BNEF"
$about.Range("B9").Value = "Lithium-ion Battery Pack Prices Hit Record Low of `$139/kWh"
$about.Range("B10").Value = 2023
$about.Range("B10").HorizontalAlignment = -4131
$about.Range("B11").Value = "https://about.bnef.com/blog/lithium-ion-battery-pack-prices-hit-record-low-of-139-kwh/#:~:text=Given%20this%2C%20BNEF%20expects%20average,and%20%2480%2FkWh%20in%202030."

# New conversion factor row used by the updated BPP formula.
$about.Range("A17").Formula = "=1/0.951"
$about.Range("B17").Value = "2022 to 2023"

# Update the BPP sheet 2023 price using the new BNEF figure.
$bpp.Range("D2").Formula = "=139/(About!A16*About!A17)"

# Refresh sheet selections to match the saved workbook state.
$bpp.Activate()
$bpp.Range("D3").Select()

$about.Activate()
$about.Range("A12:XFD12").Select()

$wb.Save()
